$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 47

# Column A: date serial number, formatted like the cells above it (style from A46)
$ws.Cells.Item($newRow, 1).Value = 45951
$ws.Range("A46").Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Columns B-E: text values (quotes with comma decimal separator), stored as text
$ws.Cells.Item($newRow, 2).Value = "21,7178"
$ws.Cells.Item($newRow, 3).Value = "15,6198"
$ws.Cells.Item($newRow, 4).Value = "15,4273"
$ws.Cells.Item($newRow, 5).Value = "15,4273"
